$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price column as Text so numeric-looking strings (e.g. "1.00", "568.97")
# are stored as literal text instead of being parsed into numbers, matching the
# original workbook where every Price/Volume cell is an inline string.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "64.858.33"
$ws.Range("E2").Value = "  -2.53%  "

# Row 3
$ws.Range("D3").Value = "3.167.13"
$ws.Range("E3").Value = "  -7.55%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "568.97"
$ws.Range("E5").Value = "  -2.82%  "

# Row 6
$ws.Range("D6").Value = "170.52"
$ws.Range("E6").Value = "  -6.11%  "

# Row 7
$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  -0.71%  "

# Row 8
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("D9").Value = "3.168.08"
$ws.Range("E9").Value = "  -7.46%  "

# Row 10
$ws.Range("D10").Value = "0.124"
$ws.Range("E10").Value = "  -5.97%  "

# Row 11
$ws.Range("D11").Value = "6.56"
$ws.Range("E11").Value = "  -6.08%  "

# Row 12
$ws.Range("D12").Value = "0.395"
$ws.Range("E12").Value = "  -4.65%  "

# Row 13
$ws.Range("D13").Value = "3.715.98"
$ws.Range("E13").Value = "  -7.57%  "

# Row 14
$ws.Range("E14").Value = "  +1.15%  "

# Row 15
$ws.Range("D15").Value = "27.15"
$ws.Range("E15").Value = "  -6.93%  "

# Row 16
$ws.Range("D16").Value = "64.878.86"
$ws.Range("E16").Value = "  -2.35%  "

# Row 17
$ws.Range("D17").Value = "0.0000162"
$ws.Range("E17").Value = "  -6.09%  "

# Row 18
$ws.Range("D18").Value = "3.166.68"
$ws.Range("E18").Value = "  -7.44%  "

# Row 19
$ws.Range("E19").Value = "  -3.09%  "

# Row 20
$ws.Range("D20").Value = "12.85"
$ws.Range("E20").Value = "  -7.24%  "

# Row 21
$ws.Range("D21").Value = "356.69"
$ws.Range("E21").Value = "  -3.41%  "

# Row 22
$ws.Range("E22").Value = "  -4.43%  "

# Row 23
$ws.Range("E23").Value = "  +0.38%  "

# Row 24
$ws.Range("D24").Value = "69.26"
$ws.Range("E24").Value = "  -5.40%  "

# Row 25
$ws.Range("D25").Value = "0.498"
$ws.Range("E25").Value = "  -6.74%  "

# Row 26
$ws.Range("E26").Value = "  -7.69%  "

# Row 27
$ws.Range("D27").Value = "9.69"
$ws.Range("E27").Value = "  -1.47%  "

# Row 28
$ws.Range("E28").Value = "  -2.24%  "

# Row 29
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.10%  "

# Row 30
$ws.Range("E30").Value = "  -0.18%  "

# Row 31
$ws.Range("D31").Value = "1.91"
$ws.Range("E31").Value = "  -4.33%  "

# Row 32
$ws.Range("D32").Value = "5.36"
$ws.Range("E32").Value = "  -7.76%  "

# Row 33
$ws.Range("D33").Value = "21.98"
$ws.Range("E33").Value = "  -6.03%  "

# Row 34
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "1.21"
$ws.Range("E34").Value = "  -5.12%  "

# Row 35
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "6.65"
$ws.Range("E35").Value = "  -6.13%  "

# Row 36
$ws.Range("E36").Value = "  -6.87%  "

# Row 37
$ws.Range("D37").Value = "154.81"
$ws.Range("E37").Value = "  -5.09%  "

# Row 38
$ws.Range("D38").Value = "0.838"
$ws.Range("E38").Value = "  -3.38%  "

# Row 39
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "26.19"
$ws.Range("E39").Value = "  -5.03%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "1.76"
$ws.Range("E40").Value = "  -2.49%  "

# Row 41
$ws.Range("D41").Value = "2.50"
$ws.Range("E41").Value = "  -5.66%  "

# Row 42
$ws.Range("D42").Value = "2.663.25"
$ws.Range("E42").Value = "  -2.08%  "

# Row 43
$ws.Range("D43").Value = "4.20"
$ws.Range("E43").Value = "  -5.31%  "

# Row 44
$ws.Range("E44").Value = "  -4.23%  "

# Row 45
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "39.42"
$ws.Range("E45").Value = "  -1.25%  "

# Row 46
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "0.0659"
$ws.Range("E46").Value = "  -4.52%  "

# Row 47
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "24.24"
$ws.Range("E47").Value = "  -3.23%  "

# Row 48
$ws.Range("D48").Value = "323.91"
$ws.Range("E48").Value = "  -3.53%  "

# Row 49
$ws.Range("D49").Value = "0.0274"
$ws.Range("E49").Value = "  -4.65%  "

# Row 50
$ws.Range("E50").Value = "  -1.56%  "

# Row 51
$ws.Range("E51").Value = "  +0.02%  "

# Restore the default style on the Price column so cells that were force-formatted
# as Text above end up with the same (default) style as the rest of the workbook.
$ws.Range("D2:D51").Style = "Normal"
